# Minor PV cost updates
$wb = $excel.ActiveWorkbook

# --- ScenarioParameters sheet ---
$wsParams = $wb.Worksheets.Item("ScenarioParameters")

# PV cost adjustment factor tweak on row 3 (G3)
$wsParams.Range("G3").Value = 0.095

# New PV_Cost_adjust value for row 4 (H4) that didn't previously exist
$wsParams.Range("H4").Value = 0.75

# Move the cursor/selection on this sheet and leave it as the non-active tab
$wsParams.Activate()
$wsParams.Range("H7").Select()

# --- SpecsDataCalib sheet ---
$wsCalib = $wb.Worksheets.Item("SpecsDataCalib")

# Updated PV cost figures (I2: PopEndYearHigh-ish / J2 column) following recalculation
$wsCalib.Range("I2").Value = 26858617.899999999
$wsCalib.Range("J2").Value = 26026616.100000001

# This sheet becomes the active tab, with the cursor resting on J2
$wsCalib.Activate()
$wsCalib.Range("J2").Select()
